$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.608.63'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '3.088.91'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.086.51'
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("E13").Value = '  -4.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.122'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").Value = '3.604.64'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '66.583.27'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.91%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.34%  '
$ws.Range("D20").Value = '3.090.81'
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '483.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.686'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.24%  '
$ws.Range("E26").Value = '  -3.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  -4.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("E33").Value = '  -2.94%  '
$ws.Range("D34").Value = '0.0₃0935'
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.941'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.309'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.34%  '
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("E42").Value = '  -5.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.84%  '
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").Value = '2.777.76'
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0345'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '367.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.06%  '
